# GoTyme interview case study - xgb threshold tuning and payoff matrix
# Adds a "Person A" threshold/confusion-matrix & payoff section (rows 68-103)
# to the "Case study" sheet, introduces a two-decimal Percent number format
# for the new ratio cells, and tidies the B43:B49 formatting so it matches
# the plain "quote-prefix" style used elsewhere on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Case study")
$ws2 = $wb.Worksheets.Item("Data Dictionary")

# ---------------------------------------------------------------------
# 1) Re-point B43:B49 at the same formatting as B23 (quote-prefix, no
#    bold) instead of the bold/quote-prefix combo they used before -
#    this frees up that style slot for the new Percent format below.
# ---------------------------------------------------------------------
$ws.Range("B23").Copy()
$ws.Range("B43:B49").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Person A - confusion matrix / threshold block (rows 68-79)
# ---------------------------------------------------------------------
$ws.Range("D68").Value = "Person A"
$ws.Range("E68").Value = 0.1
$ws.Range("F68").Value = 0.25
$ws.Range("G68").Value = 0.65

$ws.Range("B76").Value = "Actual"
$ws.Range("E70").Value = "Predicted"
$ws.Range("H76").Value = "Payoffs"

$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 1

$ws.Range("C72").Value = 0
$ws.Range("D72").Value = 5553
$ws.Range("E72").Value = 659
$ws.Range("G72").Formula = '=SUM(D72:E73)'

$ws.Range("C73").Value = 1
$ws.Range("D73").Value = 137
$ws.Range("E73").Value = 649

$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 1
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 1

$ws.Range("C76").Value = 0
$ws.Range("D76").Formula = '=D72/$G$72'
$ws.Range("E76").Formula = '=E72/$G$72'
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 300
$ws.Range("K76").Value = -300

$ws.Range("C77").Value = 1
$ws.Range("D77").Formula = '=D73/$G$72'
$ws.Range("E77").Formula = '=E73/$G$72'
$ws.Range("I77").Value = 1
$ws.Range("J77").Formula = '=0.1*-285 + 0.25*-705 + 0.65*-1225'
$ws.Range("K77").Formula = '=0.1*285 + 0.25*705 + 0.65*1225'

$ws.Range("D79").Formula = '=D76*J76+D77*J77+E76*K76+E77*K77'

# ---------------------------------------------------------------------
# 3) Second block (rows 82-91) - predictions all land in the "0" bucket
# ---------------------------------------------------------------------
$ws.Range("E82").Value = "Predicted"

$ws.Range("D83").Value = 0
$ws.Range("E83").Value = 1

$ws.Range("C84").Value = 0
$ws.Range("D84").Value = 24848
$ws.Range("E84").Value = 0
$ws.Range("G84").Formula = '=SUM(D84:E85)'

$ws.Range("C85").Value = 1
$ws.Range("D85").Value = 3143
$ws.Range("E85").Value = 0

$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 1
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 1

$ws.Range("B88").Value = "Actual"
$ws.Range("C88").Value = 0
$ws.Range("D88").Formula = '=D84/$G$84'
$ws.Range("E88").Formula = '=E84/$G$84'
$ws.Range("H88").Value = "Payoffs"
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 300
$ws.Range("K88").Value = -300

$ws.Range("C89").Value = 1
$ws.Range("D89").Formula = '=D85/$G$84'
$ws.Range("E89").Formula = '=E85/$G$84'
$ws.Range("I89").Value = 1
$ws.Range("J89").Formula = '=0.1*-285 + 0.25*-705 + 0.65*-1225'
$ws.Range("K89").Formula = '=0.1*285 + 0.25*705 + 0.65*1225'

$ws.Range("D91").Formula = '=D88*J88+D89*J89+E88*K88+E89*K89'

# ---------------------------------------------------------------------
# 4) Third block (rows 94-103) - predictions all land in the "1" bucket
# ---------------------------------------------------------------------
$ws.Range("E94").Value = "Predicted"

$ws.Range("D95").Value = 0
$ws.Range("E95").Value = 1

$ws.Range("C96").Value = 0
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 24848
$ws.Range("G96").Formula = '=SUM(D96:E97)'

$ws.Range("C97").Value = 1
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 3143

$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 1
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1

$ws.Range("B100").Value = "Actual"
$ws.Range("C100").Value = 0
$ws.Range("D100").Formula = '=D96/$G$96'
$ws.Range("E100").Formula = '=E96/$G$96'
$ws.Range("H100").Value = "Payoffs"
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 300
$ws.Range("K100").Value = -300

$ws.Range("C101").Value = 1
$ws.Range("D101").Formula = '=D97/$G$96'
$ws.Range("E101").Formula = '=E97/$G$96'
$ws.Range("I101").Value = 1
$ws.Range("J101").Formula = '=0.1*-285 + 0.25*-705 + 0.65*-1225'
$ws.Range("K101").Formula = '=0.1*285 + 0.25*705 + 0.65*1225'

$ws.Range("D103").Formula = '=D100*J100+D101*J101+E100*K100+E101*K101'

# ---------------------------------------------------------------------
# 5) Apply the two-decimal Percent number format to the ratio cells
# ---------------------------------------------------------------------
$pctRange = $ws.Range("D76:E77,D88:E89,D100:E101")
$pctRange.Style = "Percent"
$pctRange.NumberFormat = "0.00%"

# ---------------------------------------------------------------------
# 6) Selection / view bookkeeping to match the saved workbook state
# ---------------------------------------------------------------------
$ws.Range("S90").Select()
$ws2.Range("C15").Select()
$ws.Activate()
